$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Append a new log entry row (row 11), mirroring the style/shape of the
# existing rows (6-10): date in col A, hours in col B, description in col C.
$ws.Range("A11").Value = "26.3.2019"
$ws.Range("B11").Value = 3
$ws.Range("C11").Value = "TypeScript opettelua, koodin refaktorointia käyttämään paremmin TypeScriptiä, storen refaktorointia"

# Copy formatting from the row above so the new row matches the table style.
$ws.Range("A10:C10").Copy()
$ws.Range("A11:C11").PasteSpecial(-4122) # xlPasteFormats

# The new entry's description is long, so the row grows taller to fit the
# wrapped text (matches the autofit height Excel would apply on entry).
$ws.Rows.Item(11).AutoFit()

# Update the active selection to the newly added cell, matching the saved
# workbook state after the edit.
$ws.Range("C11").Select()
